$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ E=3; G=149.656361;          H=448.969083;         I=0.5921360794347563; J=0.5921360794347564; K=3; M=1.287649333333333; N=3.862948;  O=0.1870471291026542; P=0.1870471291026542; Q=192.7049134707426; R=1734.344221236684; S=0.1107573536963724;  T=0.1107573536963724  }
    3  = @{ E=3; G=149.656361;          H=448.969083;         I=0.5921360794347563; J=0.5921360794347564; K=3; M=4.622908;          N=13.868724; O=0.6715350578151914; P=0.6715350578151914; Q=691.8475885177878; R=6226.628296660091; S=0.3976401363376799;  T=0.3976401363376799  }
    4  = @{ E=3; G=149.656361;          H=448.969083;         I=0.5921360794347563; J=0.5921360794347564; K=3; M=0.9735329999999999; N=2.920599;  O=0.1414178130821545; P=0.1414178130821545; Q=145.695406093413;  R=1311.258654840717; S=0.08373858940070415; T=0.08373858940070417 }
    5  = @{ E=3; G=52.73412466666667;   H=158.202374;         I=0.208649853730866;  J=0.208649853730866;  K=3; M=1.287649333333333; N=3.862948;  O=0.1870471291026542; P=0.1870471291026542; Q=67.90306047095022;  R=611.127544238552;  S=0.0390273561280472;  T=0.0390273561280472  }
    6  = @{ E=3; G=52.73412466666667;   H=158.202374;         I=0.208649853730866;  J=0.208649853730866;  K=3; M=4.622908;          N=13.868724; O=0.6715350578151914; P=0.6715350578151914; Q=243.7850067945307;  R=2194.065061150776; S=0.1401156915882883;  T=0.1401156915882883  }
    7  = @{ E=3; G=52.73412466666667;   H=158.202374;         I=0.208649853730866;  J=0.208649853730866;  K=3; M=0.9735329999999999; N=2.920599;  O=0.1414178130821545; P=0.1414178130821545; Q=51.33841058911399;  R=462.0456953020259; S=0.02950680601453048; T=0.02950680601453048 }
    8  = @{ E=3; G=50.34932566666667;   H=151.047977;         I=0.1992140668343777; J=0.1992140668343777; K=3; M=1.287649333333333; N=3.862948;  O=0.1870471291026542; P=0.1870471291026542; Q=64.83227562846622;  R=583.490480656196;  S=0.03726241927823461; T=0.03726241927823462 }
    9  = @{ E=3; G=50.34932566666667;   H=151.047977;         I=0.1992140668343777; J=0.1992140668343777; K=3; M=4.622908;          N=13.868724; O=0.6715350578151914; P=0.6715350578151914; Q=232.7603004190387;  R=2094.842703771348; S=0.1337792298892232;  T=0.1337792298892232  }
    10 = @{ E=3; G=50.34932566666667;   H=151.047977;         I=0.1992140668343777; J=0.1992140668343777; K=3; M=0.9735329999999999; N=2.920599;  O=0.1414178130821545; P=0.1414178130821545; Q=49.01673006424699;  R=441.150570578223;  S=0.02817241766691985; T=0.02817241766691985 }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    foreach ($col in $row.Keys) {
        $ws.Range("$col$r").Value = $row[$col]
    }
}
